# Update the "2019" personal stats worksheet for Day 1 results.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2019")
$ws.Activate()

# Row 5 corresponds to Day 1. Update the puzzle title and the recorded times.
$ws.Range("B5").Value = "Day 1: The Tyranny of the Rocket Equation"

# Part 1 time (C5), My Time (E5), Leaderboard fastest (F5), expressed as day fractions (h:mm:ss)
$ws.Range("C5").Value = 142 / 86400
$ws.Range("E5").Value = 308 / 86400
$ws.Range("F5").Value = 238 / 86400

# Finish rank
$ws.Range("H5").Value = "2nd"

# Leave the active cell where the author left it after the edit.
$ws.Range("H6").Select()
